$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing PriceChange/UpDown values for the existing last row (row 3)
$ws.Range("X3").Value = -0.06999999999999318
$ws.Range("Y3").Value = "Down"

# Append a brand new data row (row 4) with the latest day's data
$ws.Range("A4").Value = 42641.890057870369
$ws.Range("B4").Value = -11
$ws.Range("C4").Value = "Sell"
$ws.Range("D4").Value = -6
$ws.Range("E4").Value = 15234
$ws.Range("F4").Value = 2341
$ws.Range("G4").Value = 55
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 41
$ws.Range("J4").Value = 58
$ws.Range("K4").Value = 13699
$ws.Range("L4").Value = 275
$ws.Range("M4").Value = 204
$ws.Range("N4").Value = 12
$ws.Range("O4").Value = 17
$ws.Range("P4").Value = "Noun"
$ws.Range("Q4").Value = 64.72814683513376
$ws.Range("R4").Value = -32.1

$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("S4").Value = -0.0755

$ws.Range("T4").NumberFormat = "0.00%"
$ws.Range("T4").Value = -0.0025

$ws.Range("U4").Value = 6.79
$ws.Range("V4").Value = 1.88
$ws.Range("W4").Value = 0
